# mod 10 & lab 2 completed / completed
#
# Restructure the worksheet:
#   - D/E header labels: "..._per_thousand_pop" -> "..._per_tenThousand_pop"
#     (the underlying formulas/values for D & E are unchanged)
#   - column H used to hold the raw avg_annual_pay figure; it now becomes
#     a derived "avg_annual_pay_thousands" column computed as =O/1000
#   - a new column O ("avg_annual_pay") holds the raw figure that used
#     to live in H
#   - sheet view bookkeeping (selection / zoom) follows the column insert

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 49

# 1) Move the raw avg_annual_pay values that currently live in column H
#    over to the new column O.
for ($r = 2; $r -le $lastRow; $r++) {
    $rawPay = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 15).Value2 = $rawPay
}

# 2) Turn H into =O/1000 (row 2 stands alone, rows 3:49 share one formula,
#    matching how the other derived columns on this sheet are laid out).
$ws.Range("H2").Formula = "=O2/1000"
$ws.Range("H3:H49").Formula = "=O3/1000"

# 3) Header row text.
$ws.Cells.Item(1, 4).Value2  = "num_employees_per_tenThousand_pop"
$ws.Cells.Item(1, 5).Value2  = "num_businesses_per_tenThousand_pop"
$ws.Cells.Item(1, 8).Value2  = "avg_annual_pay_thousands"
$ws.Cells.Item(1, 15).Value2 = "avg_annual_pay"

# 4) Sheet view bookkeeping to match the edited file (selection moves to
#    E1, view is re-zoomed to 100%).
$ws.Range("E1").Select()
$excel.ActiveWindow.Zoom = 100
